# Auto-generated script applying the diff from the commit
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 489.42856
$ws.Range("I18").Value = 237.66667
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 237.66667
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 46.33332999999999
$ws.Range("N18").Value = -2568
$ws.Range("H32").Value = 8033.533
$ws.Range("I32").Value = 499.6
$ws.Range("J32").Value = 11800.5
$ws.Range("K32").Value = 499.6
$ws.Range("L32").Value = 11800.5
$ws.Range("M32").Value = -173.6
$ws.Range("N32").Value = -12452.5
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
$ws.Range("H86").Value = 2238.9167
$ws.Range("J86").Value = 1347.3334
$ws.Range("L86").Value = 1347.3334
$ws.Range("N86").Value = -3593.3334
$ws.Range("H89").Value = 2238.9167
$ws.Range("J89").Value = 1347.3334
$ws.Range("L89").Value = 6736.666999999999
$ws.Range("N89").Value = -17968.667
$ws.Range("H113").Value = 25002176
$ws.Range("I113").Value = 28573700
$ws.Range("J113").Value = 1506
$ws.Range("K113").Value = 28573700
$ws.Range("L113").Value = 1506
$ws.Range("M113").Value = -28570446
$ws.Range("N113").Value = -8014
$ws.Range("H132").Value = 6951436.5
$ws.Range("I132").Value = 7941029.5
$ws.Range("K132").Value = 23823088.5
$ws.Range("M132").Value = -23820558.5
$ws.Range("H137").Value = 2621.2827
$ws.Range("J137").Value = 2960.52
$ws.Range("L137").Value = 8881.559999999999
$ws.Range("N137").Value = -13981.56
$ws.Range("H138").Value = 2536.3333
$ws.Range("I138").Value = 2156.7693
$ws.Range("J138").Value = 2598.0125
$ws.Range("K138").Value = 6470.3079
$ws.Range("L138").Value = 7794.037499999999
$ws.Range("M138").Value = -1330.3079
$ws.Range("N138").Value = -18074.0375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9412.525
$ws.Range("I32").Value = 6711.5557
$ws.Range("J32").Value = 21566.889
$ws.Range("K32").Value = 6711.5557
$ws.Range("L32").Value = 21566.889
$ws.Range("M32").Value = -6424.5557
$ws.Range("N32").Value = -22140.889
$ws.Range("H61").Value = 66668250
$ws.Range("I61").Value = 83334560
$ws.Range("J61").Value = 2999.6667
$ws.Range("K61").Value = 83334560
$ws.Range("L61").Value = 2999.6667
$ws.Range("M61").Value = -83334348
$ws.Range("N61").Value = -3423.6667
$ws.Range("H63").Value = 33335576
$ws.Range("I63").Value = 2295.96
$ws.Range("J63").Value = 200001980
$ws.Range("K63").Value = 2295.96
$ws.Range("L63").Value = 200001980
$ws.Range("M63").Value = -1609.96
$ws.Range("N63").Value = -200003352
$ws.Range("H66").Value = 33335576
$ws.Range("I66").Value = 2295.96
$ws.Range("J66").Value = 200001980
$ws.Range("K66").Value = 11479.8
$ws.Range("L66").Value = 1000009900
$ws.Range("M66").Value = -8047.799999999999
$ws.Range("N66").Value = -1000016764
$ws.Range("H74").Value = 1485.625
$ws.Range("I74").Value = 1051.6842
$ws.Range("K74").Value = 1051.6842
$ws.Range("M74").Value = -177.6841999999999
$ws.Range("H77").Value = 1485.625
$ws.Range("I77").Value = 1051.6842
$ws.Range("K77").Value = 5258.420999999999
$ws.Range("M77").Value = -890.4209999999994
$ws.Range("H122").Value = 1616.2858
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 3157
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 9471
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -14371
$ws.Range("H136").Value = 66668250
$ws.Range("I136").Value = 83334560
$ws.Range("J136").Value = 2999.6667
$ws.Range("K136").Value = 250003680
$ws.Range("L136").Value = 8999.000100000001
$ws.Range("M136").Value = -250001130
$ws.Range("N136").Value = -14099.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4060.7812
$ws.Range("I134").Value = 939.8570999999999
$ws.Range("J134").Value = 25907.25
$ws.Range("K134").Value = 2819.5713
$ws.Range("L134").Value = 77721.75
$ws.Range("M134").Value = -284.5712999999996
$ws.Range("N134").Value = -82791.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 90910770
$ws.Range("I16").Value = 90910770
$ws.Range("K16").Value = 90910770
$ws.Range("M16").Value = -90910483
$ws.Range("H22").Value = 78035.55499999999
$ws.Range("I22").Value = 242.5
$ws.Range("J22").Value = 140270
$ws.Range("K22").Value = 242.5
$ws.Range("L22").Value = 140270
$ws.Range("M22").Value = 107.5
$ws.Range("N22").Value = -140970
$ws.Range("H31").Value = 1675.0217
$ws.Range("I31").Value = 1528.2222
$ws.Range("J31").Value = 2203.5
$ws.Range("K31").Value = 1528.2222
$ws.Range("L31").Value = 2203.5
$ws.Range("M31").Value = -1233.2222
$ws.Range("N31").Value = -2793.5
$ws.Range("H34").Value = 1675.0217
$ws.Range("I34").Value = 1528.2222
$ws.Range("J34").Value = 2203.5
$ws.Range("K34").Value = 1528.2222
$ws.Range("L34").Value = 2203.5
$ws.Range("M34").Value = -1326.2222
$ws.Range("N34").Value = -2607.5
$ws.Range("H109").Value = 10633.667
$ws.Range("J109").Value = 10633.667
$ws.Range("L109").Value = 10633.667
$ws.Range("N109").Value = -12713.667
$ws.Range("H113").Value = 90910770
$ws.Range("I113").Value = 90910770
$ws.Range("K113").Value = 90910770
$ws.Range("M113").Value = -90908600
$ws.Range("H134").Value = 17242892
$ws.Range("I134").Value = 1523.2858
$ws.Range("K134").Value = 4569.857400000001
$ws.Range("M134").Value = -2034.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 809.8182
$ws.Range("I68").Value = 1231.6666
$ws.Range("J68").Value = 651.625
$ws.Range("K68").Value = 3694.9998
$ws.Range("L68").Value = 1954.875
$ws.Range("M68").Value = -2883.9998
$ws.Range("N68").Value = -3576.875
$ws.Range("H71").Value = 809.8182
$ws.Range("I71").Value = 1231.6666
$ws.Range("J71").Value = 651.625
$ws.Range("K71").Value = 11084.9994
$ws.Range("L71").Value = 5864.625
$ws.Range("M71").Value = -7028.999400000001
$ws.Range("N71").Value = -13976.625
$ws.Range("H96").Value = 8662.5
$ws.Range("J96").Value = 8662.5
$ws.Range("L96").Value = 25987.5
$ws.Range("N96").Value = -30105.5
$ws.Range("H129").Value = 24510740
$ws.Range("I129").Value = 41667104
$ws.Range("J129").Value = 9260639
$ws.Range("K129").Value = 125001312
$ws.Range("L129").Value = 27781917
$ws.Range("M129").Value = -124996312
$ws.Range("N129").Value = -27791917

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37503696
$ws.Range("I70").Value = 27781640
$ws.Range("J70").Value = 66669868
$ws.Range("K70").Value = 27781640
$ws.Range("L70").Value = 66669868
$ws.Range("M70").Value = -27781370
$ws.Range("N70").Value = -66670408
$ws.Range("H73").Value = 37503696
$ws.Range("I73").Value = 27781640
$ws.Range("J73").Value = 66669868
$ws.Range("K73").Value = 27781640
$ws.Range("L73").Value = 66669868
$ws.Range("M73").Value = -27780704
$ws.Range("N73").Value = -66671740
$ws.Range("H80").Value = 3579.9
$ws.Range("I80").Value = 1800
$ws.Range("K80").Value = 1800
$ws.Range("M80").Value = -802
$ws.Range("H83").Value = 3579.9
$ws.Range("I83").Value = 1800
$ws.Range("K83").Value = 9000
$ws.Range("M83").Value = -4008
$ws.Range("H102").Value = 3777.5557
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3777.5557
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 3777.5557
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -7021.5557
$ws.Range("H126").Value = 2176
$ws.Range("I126").Value = 1851.4286
$ws.Range("J126").Value = 2933.3333
$ws.Range("K126").Value = 5554.2858
$ws.Range("L126").Value = 8799.999899999999
$ws.Range("M126").Value = -3084.2858
$ws.Range("N126").Value = -13739.9999
$ws.Range("H135").Value = 38313.332
$ws.Range("J135").Value = 36852.5
$ws.Range("L135").Value = 36852.5
$ws.Range("N135").Value = -46992.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 200
$ws.Range("K23").Value = 200
$ws.Range("M23").Value = 30
$ws.Range("H40").Value = 8003
$ws.Range("I40").Value = 1004
$ws.Range("J40").Value = 11502.5
$ws.Range("K40").Value = 1004
$ws.Range("L40").Value = 11502.5
$ws.Range("M40").Value = -868
$ws.Range("N40").Value = -11774.5
$ws.Range("H54").Value = 10084
$ws.Range("J54").Value = 10084
$ws.Range("L54").Value = 10084
$ws.Range("N54").Value = -11372
$ws.Range("H55").Value = 309.7857
$ws.Range("I55").Value = 244
$ws.Range("K55").Value = 244
$ws.Range("M55").Value = -71
$ws.Range("H122").Value = 70834340
$ws.Range("I122").Value = 141666670
$ws.Range("J122").Value = 2002.5
$ws.Range("K122").Value = 425000010
$ws.Range("L122").Value = 6007.5
$ws.Range("M122").Value = -424997560
$ws.Range("N122").Value = -10907.5
$ws.Range("H136").Value = 1535.35
$ws.Range("I136").Value = 1373.7222
$ws.Range("K136").Value = 4121.1666
$ws.Range("M136").Value = -1571.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 12069
$ws.Range("J47").Value = 12069
$ws.Range("L47").Value = 12069
$ws.Range("N47").Value = -13213
$ws.Range("H48").Value = 10043.333
$ws.Range("H94").Value = 15000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16802
$ws.Range("H122").Value = 6946550
$ws.Range("I122").Value = 10871768
$ws.Range("J122").Value = 1933.4615
$ws.Range("K122").Value = 32615304
$ws.Range("L122").Value = 5800.3845
$ws.Range("M122").Value = -32612854
$ws.Range("N122").Value = -10700.3845
$ws.Range("H126").Value = 69445544
$ws.Range("I126").Value = 101010900
$ws.Range("J126").Value = 1774
$ws.Range("K126").Value = 303032700
$ws.Range("L126").Value = 5322
$ws.Range("M126").Value = -303030230
$ws.Range("N126").Value = -10262
$ws.Range("H136").Value = 1994.4667
$ws.Range("I136").Value = 1734.6666
$ws.Range("K136").Value = 5203.9998
$ws.Range("M136").Value = -2653.9998
